# Stage 1: update companies data
# The underlying rows 3-11 get re-sorted (company records shuffled into a
# new order) and a handful of SIC-code / description cells are updated to
# match. Only columns A (Company Name), B (Company Number), H (Category),
# I (SIC Codes), J (SIC Description) and K (Typical Use Case) actually
# change value; C/D/E/F/G are identical between the swapped rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Company number (B) and purely-numeric SIC codes (I) must stay text, or
# Excel's type-inference will coerce the digit string into a Number -
# force a text number format before writing those cells that are actually
# being touched (so we don't perturb the style of untouched cells).
foreach ($row in @(3, 4, 5, 6, 7, 9, 10, 11)) {
    $ws.Cells.Item($row, 2).NumberFormat = "@"
}
foreach ($row in @(4, 5, 6, 11)) {
    $ws.Cells.Item($row, 9).NumberFormat = "@"
}

$data = @{
    3  = @{ A = "GANDER INVESTMENTS LTD"; B = "16473515"; H = "Investments"; I = "68100,68209"; J = ""; K = "" }
    4  = @{ A = "SEVEN (HOLDCO) LIMITED"; B = "16473606"; H = "Other"; I = "64209"; J = "Activities of other holding companies n.e.c."; K = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles." }
    5  = @{ A = "BRIDGEWICK PARTNERS LIMITED"; B = "16473142"; H = "Partners"; I = "64999"; J = "Financial intermediation not elsewhere classified"; K = "Catch-all credit-oriented SPVs for novel lending structures." }
    6  = @{ A = "MARMIMI HOLDING LIMITED"; B = "16473234"; H = "Other"; I = "64209"; J = "Activities of other holding companies n.e.c."; K = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles." }
    7  = @{ A = "THE DISLEY GROUP LTD"; B = "16473398"; H = "Other" }
    9  = @{ A = "TLJ INVESTMENT LTD"; B = "16473151"; H = "Investments"; I = "41100,55100,68100"; J = ""; K = "" }
    10 = @{ A = "INTERCONTINENTAL HOLDING COMPANY LIMITED"; B = "16473418" }
    11 = @{ A = "GAUNT CAPITAL LTD"; B = "16473262"; H = "Capital"; I = "64209"; J = "Activities of other holding companies n.e.c."; K = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles." }
}

$colIndex = @{ A = 1; B = 2; H = 8; I = 9; J = 10; K = 11 }

foreach ($rowNum in $data.Keys) {
    $rowValues = $data[$rowNum]
    foreach ($colLetter in $rowValues.Keys) {
        $value = $rowValues[$colLetter]
        $col = $colIndex[$colLetter]
        $cell = $ws.Cells.Item($rowNum, $col)
        if (($colLetter -eq "J" -or $colLetter -eq "K") -and $value -eq "") {
            # Re-create an *empty text* cell (not a fully-blank one) to
            # mirror the original inlineStr-with-no-content cells.
            $cell.Value = "'"
        } else {
            $cell.Value = $value
        }
    }
}

$wb.Save()
